# The deck ships with two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" / "Office" colour scheme (only
#                             referenced by the notes master, never seen on
#                             a slide)
#   ppt/theme/theme2.xml  -> "Integral" / "Red Violet" colour scheme
#                             (referenced by the slide master, i.e. this is
#                             the palette every slide actually renders with)
#
# The authored edit swaps the two themes' colour schemes so the deck that
# was pink/violet ("Integral") now uses the plain default Office palette,
# while the (invisible, notes-only) spare theme keeps the old pink/violet
# values.
#
# PowerPoint's object model doesn't expose the raw theme XML parts
# directly, but it does expose the 12 theme colour slots (dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink) through ThemeColorScheme, which is
# backed by the theme that the slide master (and therefore every slide)
# actually uses. Re-pointing each slot to the "Office Theme" values
# reproduces the visible effect of the authored swap.

$p = $ppt.ActivePresentation

# Target values = the stock "Office Theme" colour scheme (what theme1.xml
# held before the edit, and what theme2.xml holds afterwards), expressed
# as COM RGB() integers (0x00BBGGRR).
$officeThemeColors = @(
    0,        # 1  dk1      000000
    16777215, # 2  lt1      FFFFFF
    6968388,  # 3  dk2      44546A
    15132391, # 4  lt2      E7E6E6
    13998939, # 5  accent1  5B9BD5
    3243501,  # 6  accent2  ED7D31
    10855845, # 7  accent3  A5A5A5
    49407,    # 8  accent4  FFC000
    12874308, # 9  accent5  4472C4
    4697456,  # 10 accent6  70AD47
    12673797, # 11 hlink    0563C1
    7491477   # 12 folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}
